{"js": "// Insert a \"Meta description\" paragraph right after the title (Heading 1),\n// and move the old title/description pair that used to sit at the very end\n// of the document: drop the duplicated bold title paragraph there, and turn\n// the now-orphaned italic paragraph into the new AI-image-prompt paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titleText = \"Play #90 Spins for Free - The Ultimate Football-themed Slot\";\nconst oldMetaText =\n  \"Celebrate the national Italian sport with Espresso Games' #90 Spins, available to play for free. Experience 6,784 ways to win and double wild reels.\";\nconst newImagePromptText =\n  \"Create a feature image that captures the essence of #90spins! The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a football and wearing a football jersey with the number 90 on it. The background should feature a football stadium filled with cheering fans waving their country's flags. Let's make the Maya warrior the hero of the game, ready to take on any opponent and score big!\";\n\nconst items = paragraphs.items;\n\n// The very first paragraph is the document title (Heading 1).\nconst titlePara = items[0];\n\n// Find the trailing duplicate bold title paragraph and the trailing italic\n// meta-description paragraph (both live at the end of the document, after\n// the \"What we don't like\" bullet list).\nlet trailingTitleIndex = -1;\nlet trailingMetaIndex = -1;\nfor (let i = items.length - 1; i >= 1; i--) {\n  const t = items[i].text;\n  if (trailingTitleIndex === -1 && t === titleText) {\n    trailingTitleIndex = i;\n  }\n  if (trailingMetaIndex === -1 && t === oldMetaText) {\n    trailingMetaIndex = i;\n  }\n  if (trailingTitleIndex !== -1 && trailingMetaIndex !== -1) break;\n}\n\nconst trailingMetaPara = items[trailingMetaIndex];\nconst trailingTitlePara = items[trailingTitleIndex];\n\n// 1) Insert the new \"Meta description: ...\" paragraph right after the title.\nconst metaPara = titlePara.insertParagraph(\n  \"Meta description\" + \": \" + oldMetaText,\n  \"After\"\n);\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Make the \"Meta description\" label bold, matching the source formatting.\nconst boldRanges = metaPara.search(\"Meta description\", { matchCase: true });\nboldRanges.load(\"items\");\nawait context.sync();\nboldRanges.items[0].font.bold = true;\nawait context.sync();\n\n// 2) Replace the trailing italic paragraph's text with the new image-prompt\n// copy, keeping its existing italic run formatting.\ntrailingMetaPara.insertText(newImagePromptText, \"Replace\");\nawait context.sync();\n\n// 3) Remove the now-redundant trailing bold title paragraph entirely.\ntrailingTitlePara.delete();\nawait context.sync();\n", "ps1": "# Insert a \"Meta description\" paragraph right after the title (Heading 1),\n# and consolidate the old title/description pair that used to sit at the very\n# end of the document: drop the duplicated bold title paragraph there, and\n# turn the now-orphaned italic paragraph into the new AI-image-prompt\n# paragraph.\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play #90 Spins for Free - The Ultimate Football-themed Slot\"\n$oldMetaText = \"Celebrate the national Italian sport with Espresso Games' #90 Spins, available to play for free. Experience 6,784 ways to win and double wild reels.\"\n$newImagePromptText = \"Create a feature image that captures the essence of #90spins! The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a football and wearing a football jersey with the number 90 on it. The background should feature a football stadium filled with cheering fans waving their country's flags. Let's make the Maya warrior the hero of the game, ready to take on any opponent and score big!\"\n\n# The very first paragraph is the document title (Heading 1).\n$titlePara = $d.Paragraphs.Item(1)\n\n# Locate the trailing duplicate bold title paragraph and the trailing italic\n# meta-description paragraph (both live at the end of the document, after\n# the \"What we don't like\" bullet list).\n$count = $d.Paragraphs.Count\n$trailingTitleIndex = -1\n$trailingMetaIndex = -1\nfor ($i = $count; $i -ge 2; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($trailingTitleIndex -eq -1 -and $t -eq $titleText) {\n        $trailingTitleIndex = $i\n    }\n    if ($trailingMetaIndex -eq -1 -and $t -eq $oldMetaText) {\n        $trailingMetaIndex = $i\n    }\n    if ($trailingTitleIndex -ne -1 -and $trailingMetaIndex -ne -1) {\n        break\n    }\n}\n\n# NOTE: the two trailing edits are performed BEFORE the new paragraph is\n# inserted near the top of the document, so the indices found above\n# (trailingTitleIndex / trailingMetaIndex) still point at the right\n# paragraphs \u2014 inserting a paragraph near the top would otherwise shift every\n# later paragraph's index by one.\n\n# 1) Replace the trailing italic paragraph's text with the new image-prompt\n# copy, keeping its existing italic run formatting. Exclude the trailing\n# paragraph mark from the range so the assignment replaces (not prepends).\n$metaRange = $d.Paragraphs.Item($trailingMetaIndex).Range\n[void]$metaRange.MoveEnd(1, -1)\n$metaRange.Text = $newImagePromptText\n\n# 2) Remove the now-redundant trailing bold title paragraph entirely.\n$d.Paragraphs.Item($trailingTitleIndex).Range.Delete()\n\n# 3) Insert the new \"Meta description: ...\" paragraph right after the title.\n[void]$titlePara.Range.InsertParagraphAfter()\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n$metaPara.Range.Text = \"Meta description\" + \": \" + $oldMetaText\n\n# Make the \"Meta description\" label bold, matching the source formatting.\n$labelRange = $metaPara.Range.Duplicate\n$labelRange.Find.ClearFormatting()\n$labelRange.Find.Text = \"Meta description\"\n[void]$labelRange.Find.Execute()\n$labelRange.Bold = 1\n"}
